$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.904526829719543
$ws.Range("B1").Value = 2.187911748886108
$ws.Range("C1").Value = 2.384790182113647
$ws.Range("D1").Value = 3.477867126464844
$ws.Range("E1").Value = 1.195002555847168
